# "deleted uranium ammo now has seperate research"
#
# Adds 4 new research entries (high/low caliber depleted-uranium shells,
# each with a base + lvl_2 upgrade) to the rebalance_localizations sheet,
# each contributing a "description" row and a "name" row (key in column A,
# localized text in column B). Also updates the scratch/staging area on
# Sheet1 to point at the newly typed key strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rebalance_localizations")
$scratch = $wb.Worksheets.Item("Sheet1")

# Last currently used row in the data table (header is row 1).
$lastRow = 761

# --- 1) Scratch area on Sheet1: user starts drafting the new keys there ---
$scratch.Cells.Item(6, 1).Value = "gui/menu/research/description/depleted_uranium_high_caliber_shells"
$scratch.Cells.Item(7, 1).Value = "gui/menu/research/name/depleted_uranium_high_caliber_shells"

# --- 2) Append the 8 new rows to the bottom of the main table ---
$rowNameHigh    = $lastRow + 1
$rowDescLow     = $lastRow + 2
$rowNameLow     = $lastRow + 3
$rowDescHigh    = $lastRow + 4
$rowDescHighL2  = $lastRow + 5
$rowNameHighL2  = $lastRow + 6
$rowDescLowL2   = $lastRow + 7
$rowNameLowL2   = $lastRow + 8

$ws.Cells.Item($rowNameHigh, 1).Value = "gui/menu/research/name/depleted_uranium_high_caliber_shells"
$ws.Cells.Item($rowNameHigh, 2).Value = "Depleted Uranium high caliber Shells"

$ws.Cells.Item($rowDescLow, 1).Value = "gui/menu/research/description/depleted_uranium_low_caliber_shells"

$ws.Cells.Item($rowNameLow, 1).Value = "gui/menu/research/name/depleted_uranium_low_caliber_shells"

$ws.Cells.Item($rowDescLow, 2).Value = "Depelted Uranium allows the manufacturing of more effective high density rounds for low caliber weapons"

$ws.Cells.Item($rowDescHigh, 1).Value = "gui/menu/research/description/depleted_uranium_high_caliber_shells"
$ws.Cells.Item($rowDescHigh, 2).Value = "Depelted Uranium allows the manufacturing of more effective high density shells for high caliber weapons"

$ws.Cells.Item($rowNameLow, 2).Value = "Depleted Uranium low caliber Rounds"

$ws.Cells.Item($rowDescHighL2, 1).Value = "gui/menu/research/description/depleted_uranium_high_caliber_shells_lvl_2"
$ws.Cells.Item($rowNameHighL2, 1).Value = "gui/menu/research/name/depleted_uranium_high_caliber_shells_lvl_2"
$ws.Cells.Item($rowDescLowL2, 1).Value = "gui/menu/research/description/depleted_uranium_low_caliber_shells-lvl_2"
$ws.Cells.Item($rowNameLowL2, 1).Value = "gui/menu/research/name/depleted_uranium_low_caliber_shells_lvl_2"

$ws.Cells.Item($rowNameHighL2, 2).Value = "Improved Depleted Uranium high caliber Shells"
$ws.Cells.Item($rowNameLowL2, 2).Value = "Improved Depleted Uranium low caliber Rounds"

$ws.Cells.Item($rowDescHighL2, 2).Value = "Improvements for the production of depleted uranium munitions that allow a more efficient production line with higher output"
$ws.Cells.Item($rowDescLowL2, 2).Value = "Improvements for the production of depleted uranium munitions that allow a more efficient production line with higher output"

# --- 3) Re-sort the whole table (A2:K769) by column A, like the author did ---
$newLastRow = $lastRow + 8
$fullRange = $ws.Range("A2:K$newLastRow")
$keyRange = $ws.Range("A2:A$newLastRow")
$fullRange.Sort($keyRange)

# --- 4) Leave the selection where the author ended up after the edit ---
$ws.Activate()
$ws.Range("A749").Select()
